# cryptos.xlsx price/volume refresh (GitHub Actions bot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13/14 swapped rank order: Litecoin <-> WrappedEther
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'96.61"
$ws.Range("E13").Value = "  +2.05%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.895.54"
$ws.Range("E14").Value = "  +0.45%  "

# Price (D) / Volume(1h) (E) refresh for all other rows
$ws.Range("D2").Value = "27.237.11"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.901.65"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'307.82"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.5208"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "'0.3772"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").Value = "'0.07281"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "'21.19"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").Value = "'0.9032"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'0.08273"
$ws.Range("E12").Value = "  +8.32%  "
$ws.Range("D15").Value = "'5.284"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "'0.000008629"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "'14.58"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "'0.9999"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "27.241.77"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "'5.095"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "2.150.28"
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").Value = "'10.67"
$ws.Range("D24").Value = "'6.426"
$ws.Range("D25").Value = "'2.322"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").Value = "'147.14"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").Value = "'1.749"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "'18.24"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").Value = "'115.31"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").Value = "'4.836"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").Value = "'4.905"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "'0.09251"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").Value = "'0.05077"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").Value = "'0.8003"
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").Value = "'3.434"
$ws.Range("E36").Value = "  +4.76%  "
$ws.Range("D37").Value = "'2.946"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "'2.603"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "'0.5718"
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("D40").Value = "'0.02002"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "'9.044"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "'6.580"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "'116.54"
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("D45").Value = "'0.1518"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "'0.4859"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "'10.16"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "'1.626"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "'37.70"
$ws.Range("D51").Value = "'63.90"
$ws.Range("E51").Value = "  -0.07%  "
